# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 05:20"

# --- Rows 58-65: Honduras reinserted right after Serbia (before Azerbaiyan),
# shifting Azerbaiyan / Corea del Sur / Dinamarca down one slot each.
# Row 58 (Ghana) stays put but its stats are refreshed.
Set-Row 58 @("Ghana", 14154, 147, 10473, 3596, 0, 0, 85)

Set-Row 61 @("Honduras", 12769, 463, 1293, 11113, 0, 5, 363)
Set-Row 62 @("Azerbaiyan", 12729, 0, 6799, 5776, 0, 0, 154)
Set-Row 63 @("Corea del Sur", 12438, 17, 10881, 1277, 0, 0, 280)
Set-Row 64 @("Dinamarca", 12391, 0, 11282, 509, 0, 0, 600)

# Row 74 (Australia) stays put but its stats are refreshed.
Set-Row 74 @("Australia", 7474, 13, 6903, 469, 0, 0, 102)

# --- Fiyi / Dominica swap position (identical stats, so only the labels move)
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# --- Groenlandia / Islas Malvinas swap position (identical stats, so only the labels move)
$ws.Range("A207").Value = "Islas Malvinas"
$ws.Range("A208").Value = "Groenlandia"

# --- Papua Nueva Guinea / Islas Virgenes Britanicas swap position (stats differ)
Set-Row 213 @("Islas Virgenes Britanicas", 8, 0, 7, 0, 0, 0, 1)
Set-Row 214 @("Papua Nueva Guinea", 8, 0, 8, 0, 0, 0, 0)
